# Update the marksheet figures in the "quiz" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: Right count 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row: total marks 39 -> 65
$ws.Range("B12").Value = 65

# "Total" row: correct/total display "34/84" -> "65/140"
$ws.Range("E12").Value = "65/140"
